# Case_4_13 line/parallel.xlsx edit: extend the table from columns A:O to A:Q
# (add two more "parallel" trials, columns P and Q), and rebalance the
# contingency counts in columns I, K, M, O for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): add P1 = 14, Q1 = 15, matching the bold/bordered
# header style already used by B1:O1 (copy format from O1, then set value).
$ws.Range("O1").Copy()
$ws.Range("P1").PasteSpecial(-4122)
$ws.Range("P1").Value = 14

$ws.Range("O1").Copy()
$ws.Range("Q1").PasteSpecial(-4122)
$ws.Range("Q1").Value = 15

# --- Data rows 2-25: swap the I/K and M/O pairs (1<->2) and populate the
# two new columns P and Q with 2.
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = 2    # I: 1 -> 2
    $ws.Cells.Item($r, 11).Value = 1   # K: 2 -> 1
    $ws.Cells.Item($r, 13).Value = 2   # M: 1 -> 2
    $ws.Cells.Item($r, 15).Value = 1   # O: 2 -> 1
    $ws.Cells.Item($r, 16).Value = 2   # P: new column
    $ws.Cells.Item($r, 17).Value = 2   # Q: new column
}
